$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the previously-empty H6 value (sanitized global work size input),
# which causes the dependent shared formula in I6 to resolve from #DIV/0! to a real value.
$ws.Range("H6").Value = 2987.895

# Update the active cell selection to match the saved state.
$ws.Range("H11").Select()
